$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 252, shifting existing rows 252-320 down to 253-321
$ws.Rows.Item(252).Insert()

# Populate the newly inserted row 252 with the new record's data
$ws.Cells.Item(252, 1).Value = 4
$ws.Cells.Item(252, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(252, 3).Value = "Los Lagos"
$ws.Cells.Item(252, 4).Value = 44855
$ws.Cells.Item(252, 4).NumberFormat = $ws.Cells.Item(253, 4).NumberFormat
$ws.Cells.Item(252, 5).Value = 10
$ws.Cells.Item(252, 6).Value = "Fruta"
$ws.Cells.Item(252, 7).Value = 100108
$ws.Cells.Item(252, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(252, 9).Value = 100108005
$ws.Cells.Item(252, 10).Value = "Piña"
$ws.Cells.Item(252, 11).Value = "Caramelo"
$ws.Cells.Item(252, 12).Value = "Primera"
$ws.Cells.Item(252, 13).Value = 240
$ws.Cells.Item(252, 14).Value = 23000
$ws.Cells.Item(252, 15).Value = 23500
$ws.Cells.Item(252, 16).Value = 23250
$ws.Cells.Item(252, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(252, 18).Value = "Ecuador"
$ws.Cells.Item(252, 19).Value = 1938
$ws.Cells.Item(252, 20).Value = 12
